$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.82"
$ws.Range("E2").Value = "'-3.73%"
$ws.Range("D3").Value = "'30.78"
$ws.Range("E3").Value = "'-6.23%"
$ws.Range("D4").Value = "'4.951"
$ws.Range("E4").Value = "'-0.29%"
$ws.Range("D5").Value = "'0.07202"
$ws.Range("E5").Value = "'-7.55%"
$ws.Range("D6").Value = "'1.797"
$ws.Range("E6").Value = "'-9.08%"
$ws.Range("D7").Value = "'7.672"
$ws.Range("E7").Value = "'-2.41%"
$ws.Range("D8").Value = "'3.741"
$ws.Range("E8").Value = "'-1.49%"
$ws.Range("D9").Value = "'0.8972"
$ws.Range("E9").Value = "'-3.18%"
$ws.Range("D10").Value = "'0.1651"
$ws.Range("E10").Value = "'-6.16%"
$ws.Range("D11").Value = "'0.07665"
$ws.Range("E11").Value = "'-2.43%"
$ws.Range("D12").Value = "'0.08070"
$ws.Range("E12").Value = "'-6.04%"
$ws.Range("D13").Value = "'0.03037"
$ws.Range("E13").Value = "'-3.42%"
$ws.Range("E14").Value = "'0.07%"
$ws.Range("D15").Value = "'0.001499"
$ws.Range("E15").Value = "'-1.70%"
$ws.Range("D16").Value = "'0.005685"
$ws.Range("E16").Value = "'-4.02%"
$ws.Range("D17").Value = "'3.466"
$ws.Range("E17").Value = "'0.09%"
$ws.Range("E18").Value = "'-3.31%"
$ws.Range("D19").Value = "'0.3312"
$ws.Range("E19").Value = "'-0.63%"
$ws.Range("E20").Value = "'0.94%"
$ws.Range("D21").Value = "'4.040"
$ws.Range("E21").Value = "'-6.77%"
$ws.Range("D22").Value = "'0.2184"
$ws.Range("E22").Value = "'9.62%"
$ws.Range("D23").Value = "'0.04507"
$ws.Range("E23").Value = "'-1.17%"
$ws.Range("E24").Value = "'-0.93%"
$ws.Range("D25").Value = "'0.004019"
$ws.Range("E25").Value = "'-9.49%"
$ws.Range("E26").Value = "'0.00%"
$ws.Range("D39").Value = "'0.01601"
$ws.Range("E39").Value = "'-6.78%"
$ws.Range("D40").Value = "'0.04397"
$ws.Range("E40").Value = "'-6.91%"
$ws.Range("D41").Value = "'0.007395"
$ws.Range("E41").Value = "'-6.23%"
$ws.Range("D42").Value = "'0.1311"
$ws.Range("E42").Value = "'-3.18%"
$ws.Range("D43").Value = "'0.007664"
$ws.Range("E44").Value = "'-12.40%"
$ws.Range("D45").Value = "'0.009210"
$ws.Range("E45").Value = "'-12.72%"
$ws.Range("D46").Value = "'0.00005907"
$ws.Range("E46").Value = "'-5.57%"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E49").Value = "'-3.30%"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E51").Value = "'-0.01%"
